$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: update in place to the new Finland match (values shifted from old row 4 data) ---
$ws.Range("A4").Value = "lrf80kOO"
$ws.Range("B4").Value = "30/10/2024"
$ws.Range("C4").Value = "13:00"
$ws.Range("D4").Value = "FINLAND - VEIKKAUSLIIGA"
$ws.Range("E4").Value = "Haka"
$ws.Range("F4").Value = "SJK"
$ws.Range("G4").Value = 3.6
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 1.95
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 2.3
$ws.Range("L4").Value = 2.6
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 4.33
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.25
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 2.2
$ws.Range("W4").Value = 13
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 13
$ws.Range("Z4").Value = 41
$ws.Range("AA4").Value = 26
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 12
$ws.Range("AF4").Value = 41
$ws.Range("AG4").Value = 151
$ws.Range("AH4").Value = 9
$ws.Range("AI4").Value = 11
$ws.Range("AJ4").Value = 8.5
$ws.Range("AK4").Value = 17
$ws.Range("AL4").Value = 15
$ws.Range("AM4").Value = 21
$ws.Range("AN4").Value = 5.5
$ws.Range("AO4").Value = 19
$ws.Range("AP4").Value = 23
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 67
$ws.Range("AS4").Value = 151
$ws.Range("AT4").Value = 3.25
$ws.Range("AU4").Value = 7.5
$ws.Range("AV4").Value = 41
$ws.Range("AW4").Value = 4.33
$ws.Range("AX4").Value = 10
$ws.Range("AY4").Value = 19
$ws.Range("AZ4").Value = 34
$ws.Range("BA4").Value = 51
$ws.Range("BB4").Value = 101
$ws.Range("BC4").Value = 451
$ws.Range("BD4").Value = 81

# --- Row 5 (new row): the original Slovakia match, with updated odds ---
$ws.Range("A5").Value = "rFf4UJrf"
$ws.Range("B5").Value = "30/10/2024"
$ws.Range("C5").Value = "12:00"
$ws.Range("D5").Value = "SLOVAKIA - NIKE LIGA"
$ws.Range("E5").Value = "Skalica"
$ws.Range("F5").Value = "Slovan Bratislava"
$ws.Range("G5").Value = 6.1
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 1.47
$ws.Range("J5").Value = 5.5
$ws.Range("K5").Value = 2.45
$ws.Range("L5").Value = 1.93
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 9.75
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 4.5
$ws.Range("Q5").Value = 1.57
$ws.Range("R5").Value = 2.37
$ws.Range("S5").Value = 1.3
$ws.Range("T5").Value = 3.35
$ws.Range("U5").Value = 1.7
$ws.Range("V5").Value = 2.05
$ws.Range("W5").Value = 18
$ws.Range("X5").Value = 45
$ws.Range("Y5").Value = 20
$ws.Range("Z5").Value = 150
$ws.Range("AA5").Value = 60
$ws.Range("AB5").Value = 50
$ws.Range("AC5").Value = 9.75
$ws.Range("AD5").Value = 9.5
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 65
$ws.Range("AG5").Value = 450
$ws.Range("AH5").Value = 8
$ws.Range("AI5").Value = 8.25
$ws.Range("AJ5").Value = 8.75
$ws.Range("AK5").Value = 11.25
$ws.Range("AL5").Value = 11.75
$ws.Range("AM5").Value = 23
$ws.Range("AN5").Value = 7.6
$ws.Range("AO5").Value = 32
$ws.Range("AP5").Value = 30
$ws.Range("AQ5").Value = 175
$ws.Range("AR5").Value = 175
$ws.Range("AS5").Value = 350
$ws.Range("AT5").Value = 3.35
$ws.Range("AU5").Value = 7.4
$ws.Range("AV5").Value = 55
$ws.Range("AW5").Value = 3.5
$ws.Range("AX5").Value = 6.7
$ws.Range("AY5").Value = 14.5
$ws.Range("AZ5").Value = 18.5
$ws.Range("BA5").Value = 40
$ws.Range("BB5").Value = 175
$ws.Range("BC5").Value = 51
$ws.Range("BD5").Value = 51
